$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '''27.155.41'
$ws.Range("E2").Value = '  -1.83%  '
$ws.Range("D3").Value = '''1.822.32'
$ws.Range("E3").Value = '  -1.28%  '
$ws.Range("E4").Value = '  -0.55%  '
$ws.Range("D5").Value = '''312.46'
$ws.Range("E5").Value = '  -1.92%  '
$ws.Range("E7").Value = '  -1.76%  '
$ws.Range("D8").Value = '''0.3692'
$ws.Range("E8").Value = '  -1.51%  '
$ws.Range("D9").Value = '''0.07241'
$ws.Range("E9").Value = '  -1.40%  '
$ws.Range("D10").Value = '''0.8571'
$ws.Range("E10").Value = '  -2.63%  '
$ws.Range("D11").Value = '''20.98'
$ws.Range("E11").Value = '  -2.65%  '
$ws.Range("D12").Value = '''1.816.79'
$ws.Range("E12").Value = '  -1.79%  '
$ws.Range("D13").Value = '''6.703'
$ws.Range("E13").Value = '  -0.37%  '
$ws.Range("D14").Value = '''0.07093'
$ws.Range("E14").Value = '  +0.00%  '
$ws.Range("D15").Value = '''5.299'
$ws.Range("E15").Value = '  -2.79%  '
$ws.Range("D16").Value = '''89.16'
$ws.Range("E16").Value = '  +1.73%  '
$ws.Range("D17").Value = '''1.006'
$ws.Range("E17").Value = '  -0.63%  '
$ws.Range("D18").Value = '''0.000008843'
$ws.Range("E18").Value = '  -1.46%  '
$ws.Range("E19").Value = '  -0.43%  '
$ws.Range("B20").Value = 'Avalanche'
$ws.Range("C20").Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range("D20").Value = '''15.03'
$ws.Range("E20").Value = '  -2.65%  '
$ws.Range("B21").Value = 'WrappedBTC'
$ws.Range("C21").Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range("D21").Value = '''27.241.28'
$ws.Range("E21").Value = '  -1.51%  '
$ws.Range("B22").Value = 'Uniswap'
$ws.Range("C22").Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range("D22").Value = '''5.126'
$ws.Range("E22").Value = '  -2.38%  '
$ws.Range("B23").Value = 'Cosmos'
$ws.Range("C23").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D23").Value = '''10.88'
$ws.Range("E23").Value = '  -2.48%  '
$ws.Range("B24").Value = 'WrappedliquidstakedEther2.0'
$ws.Range("C24").Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range("D24").Value = '''2.048.90'
$ws.Range("E24").Value = '  -1.46%  '
$ws.Range("B25").Value = 'Toncoin'
$ws.Range("C25").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D25").Value = '''1.988'
$ws.Range("E25").Value = '  -2.42%  '
$ws.Range("B26").Value = 'Monero'
$ws.Range("C26").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D26").Value = '''152.42'
$ws.Range("E26").Value = '  -1.89%  '
$ws.Range("B27").Value = 'LidoDAOToken'
$ws.Range("C27").Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range("D27").Value = '''2.187'
$ws.Range("E27").Value = '  +2.19%  '
$ws.Range("B28").Value = 'EthereumClassic'
$ws.Range("C28").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D28").Value = '''18.38'
$ws.Range("E28").Value = '  -0.79%  '
$ws.Range("B29").Value = 'InternetComputer(DFINITY)'
$ws.Range("C29").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D29").Value = '''5.221'
$ws.Range("E29").Value = '  -3.01%  '
$ws.Range("B30").Value = 'BitcoinCash'
$ws.Range("C30").Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range("D30").Value = '''116.32'
$ws.Range("E30").Value = '  -3.31%  '
$ws.Range("B31").Value = 'Stellar'
$ws.Range("C31").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D31").Value = '''0.08834'
$ws.Range("E31").Value = '  -0.82%  '
$ws.Range("B32").Value = 'ARBITRUM'
$ws.Range("C32").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D32").Value = '''1.187'
$ws.Range("E32").Value = '  -3.22%  '
$ws.Range("B33").Value = 'ImmutableX'
$ws.Range("C33").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D33").Value = '''0.7490'
$ws.Range("E33").Value = '  -3.95%  '
$ws.Range("B34").Value = 'Filecoin'
$ws.Range("C34").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D34").Value = '''4.435'
$ws.Range("E34").Value = '  -2.61%  '
$ws.Range("B35").Value = 'HuobiToken'
$ws.Range("C35").Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range("D35").Value = '''2.805'
$ws.Range("E35").Value = '  -3.48%  '
$ws.Range("B36").Value = 'Frax'
$ws.Range("C36").Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
$ws.Range("D36").Value = '''1.005'
$ws.Range("E36").Value = '  -0.45%  '
$ws.Range("B37").Value = 'TrustWalletToken'
$ws.Range("C37").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D37").Value = '''1.116'
$ws.Range("E37").Value = '  -1.86%  '
$ws.Range("B38").Value = 'VeChain'
$ws.Range("C38").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D38").Value = '''0.01968'
$ws.Range("E38").Value = '  -0.01%  '
$ws.Range("B39").Value = 'Hedera'
$ws.Range("C39").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D39").Value = '''0.05230'
$ws.Range("E39").Value = '  -1.94%  '
$ws.Range("B40").Value = 'FraxShare'
$ws.Range("C40").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D40").Value = '''7.339'
$ws.Range("E40").Value = '  +1.53%  '
$ws.Range("B41").Value = 'MXToken'
$ws.Range("C41").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range("D41").Value = '''2.877'
$ws.Range("E41").Value = '  +0.44%  '
$ws.Range("B42").Value = 'Algorand'
$ws.Range("C42").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D42").Value = '''0.1692'
$ws.Range("E42").Value = '  +0.95%  '
$ws.Range("B43").Value = 'TheSandbox'
$ws.Range("C43").Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range("D43").Value = '''0.5027'
$ws.Range("E43").Value = '  -2.58%  '
$ws.Range("B44").Value = 'Aptos'
$ws.Range("C44").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D44").Value = '''8.660'
$ws.Range("E44").Value = '  -2.96%  '
$ws.Range("B45").Value = 'EnergySwap'
$ws.Range("C45").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D45").Value = '''10.65'
$ws.Range("E45").Value = '  -0.31%  '
$ws.Range("B46").Value = 'Quant'
$ws.Range("C46").Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range("D46").Value = '''106.41'
$ws.Range("E46").Value = '  -3.48%  '
$ws.Range("B47").Value = 'Decentraland'
$ws.Range("C47").Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range("D47").Value = '''0.4733'
$ws.Range("E47").Value = '  +0.11%  '
$ws.Range("B48").Value = 'PaxDollar'
$ws.Range("C48").Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range("D48").Value = '''1.005'
$ws.Range("E48").Value = '  -0.45%  '
$ws.Range("B49").Value = 'Cronos'
$ws.Range("C49").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D49").Value = '''0.06395'
$ws.Range("E49").Value = '  -1.57%  '
$ws.Range("B50").Value = 'NEARProtocol'
$ws.Range("C50").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D50").Value = '''1.663'
$ws.Range("E50").Value = '  -2.05%  '
$ws.Range("B51").Value = 'RenderToken'
$ws.Range("C51").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D51").Value = '''1.859'
$ws.Range("E51").Value = '  -1.80%  '
